$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("evr")

# Update existing explained-variance / accumulated-variance values
$ws.Range("B2").Value = 0.22722043110768661
$ws.Range("C2").Value = 0.82989265795706235
$ws.Range("B3").Value = 0.19346937866961469
$ws.Range("B4").Value = 0.15376008888563711

# Append two new principal components (rows 5 and 6)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 0.13534970384342929

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 0.12009305545069469

# Match the style used by the rest of the data cells
$ws.Range("A4").Copy()
$ws.Range("A5:A6").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B5:B6").PasteSpecial(-4122)
